$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure every touched cell keeps its original text storage (the source workbook
# stores these as inline strings, not numbers/percentages/dates) by forcing a
# Text number format before writing the value.
$targetCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D5",
    "E5",
    "E6",
    "E7",
    "E8",
    "E9",
    "D10",
    "E10",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "E14",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "B19",
    "C19",
    "D19",
    "E19",
    "B20",
    "C20",
    "D20",
    "E20",
    "D21",
    "E21",
    "E22",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "E26",
    "E27",
    "E28",
    "D29",
    "E29",
    "D30",
    "E30",
    "E31",
    "D32",
    "E32",
    "D33",
    "E33",
    "D34",
    "E34",
    "D35",
    "E35",
    "D36",
    "E36",
    "D38",
    "E38",
    "D39",
    "E39",
    "E40",
    "D41",
    "E41",
    "E42",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51",
)
foreach ($cellRef in $targetCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.605.85"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "1.589.18"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("D5").Value = "210.58"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").Value = "19.59"
$ws.Range("E10").Value = "  -3.36%  "
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").Value = "1.811.61"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "1.577.18"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("E15").Value = "  -3.45%  "
$ws.Range("D16").Value = "64.75"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "26.604.58"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "208.45"
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.00"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "6.73"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("D24").Value = "8.86"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").Value = "146.66"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").Value = "  -2.91%  "
$ws.Range("D29").Value = "15.29"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").Value = "0.0506"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("D33").Value = "0.661"
$ws.Range("E33").Value = "  +21.16%  "
$ws.Range("D34").Value = "2.90"
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("D35").Value = "1.305.77"
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D38").Value = "0.0172"
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("D39").Value = "0.829"
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").Value = "0.789"
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("E42").Value = "  +2.76%  "
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").Value = "62.76"
$ws.Range("E44").Value = "  -3.74%  "
$ws.Range("D45").Value = "1.724.48"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").Value = "89.80"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").Value = "0.838"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("D49").Value = "0.0979"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").Value = "0.0501"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").Value = "7.53"
$ws.Range("E51").Value = "  +0.03%  "
